# Insert a new data row at row 160 (pushes the existing rows 160..217 down
# to 161..218, and expands the used range from A1:R217 to A1:R218 — matching
# the target diff), then populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(160).Insert()

$ws.Range("A160").Value = 11
$ws.Range("B160").Value = "Vega Monumental Concepción"
$ws.Range("C160").Value = "Bíobío"
$ws.Range("D160").Value = 44588
$ws.Range("E160").Value = 8
$ws.Range("F160").Value = 100112027
$ws.Range("G160").Value = "Melón"
$ws.Range("H160").Value = "Tuna"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 2200
$ws.Range("K160").Value = 600
$ws.Range("L160").Value = 700
$ws.Range("M160").Value = 655
$ws.Range("N160").Value = "$/unidad"
$ws.Range("O160").Value = "Región de O'Higgins"
$ws.Range("P160").Value = 655
$ws.Range("Q160").Value = 1
$ws.Range("R160").Value = "Hortaliza"
